$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.396.90'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.849.97'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("D4").Value = '''0.9990'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''240.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '''0.6282'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.07627'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").Value = '''0.2905'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("D10").Value = '''24.76'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").Value = '''0.07743'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '''5.034'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '''0.6800'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = '''0.00001058'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.79%  '
$ws.Range("D15").Value = '''83.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '''6.166'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '29.425.11'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '''227.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '''12.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").Value = '''0.9997'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '''7.488'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''158.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").Value = '''8.409'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").Value = '''17.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '''1.414'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.53%  '
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("D29").Value = '''0.05609'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.62%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").Value = '''4.072'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '''1.836'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").Value = '''0.7015'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").Value = '''2.589'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = '''0.01806'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.39%  '
$ws.Range("D37").Value = '1.232.19'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '''2.728'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").Value = '''6.379'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("D40").Value = '''0.9016'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.11%  '
$ws.Range("D41").Value = '''1.001'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '''101.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = '''66.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").Value = '''7.224'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.86%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.4005'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000116'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.07%  '
$ws.Range("D47").Value = '''9.022'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = '''1.679'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").Value = '''0.1139'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").Value = '''0.05706'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '''0.4629'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
